# Add clone row for zalando receipt
#
# Appends two new rows under the existing header/sample rows on Sheet1:
#   - Row 3: a "jibo" row-id cell in A3, the numeric sample (34) carried
#            over in B3, and "test" filling the remaining columns C3:S3.
#   - Row 4: a full "jibo" row spanning B4:S4 (clone marker row).
# Also updates the sheet's used-range dimension and the active selection
# to match (B4:S4, with B4 as the active cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: A3 = "jibo", B3 = 34, C3:S3 = "test" -----------------------
$ws.Range("A3").Value = "jibo"
$ws.Range("B3").Value = 34
$ws.Range("C3:S3").Value = "test"

# --- Row 4: B4:S4 = "jibo" ----------------------------------------------
$ws.Range("B4:S4").Value = "jibo"

# --- Selection / view state ---------------------------------------------
$ws.Range("B4:S4").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
